$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data for 7 Mayo (municipios not yet merged into the general CSV)
$ws.Range("A43").Value = 43957
$ws.Range("B43").Value = 61432
$ws.Range("C43").Value = 14870
$ws.Range("D43").Formula = "=B43-B42"
$ws.Range("E43").Formula = "=C43-C42"
$ws.Range("F43").Formula = "=E43/D43"

# Match the selection shown in the saved workbook (selection moves down with the new row)
$ws.Range("D42:F43").Select()
